# Weekly update: insert 3 new price rows (one week's worth of "Palta" /
# Hass / Perú data) at the top of the date-descending block that starts
# at row 831, pushing the existing rows (831-886) down to (834-889).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 831 (Excel shifts rows
# 831-886 down to 834-889, copying formatting from the row above, just
# like a manual "Insert Copied/Rows" in the UI).
$ws.Range("A831:A833").EntireRow.Insert()

# Row 831 - Especial
$ws.Cells.Item(831, 1).Value = 11
$ws.Cells.Item(831, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(831, 3).Value = "Bíobío"
$ws.Cells.Item(831, 4).Value = 45021
$ws.Cells.Item(831, 5).Value = 8
$ws.Cells.Item(831, 6).Value = "Fruta"
$ws.Cells.Item(831, 7).Value = 100106
$ws.Cells.Item(831, 8).Value = "Oleaginosos"
$ws.Cells.Item(831, 9).Value = 100106002
$ws.Cells.Item(831, 10).Value = "Palta"
$ws.Cells.Item(831, 11).Value = "Hass"
$ws.Cells.Item(831, 12).Value = "Especial"
$ws.Cells.Item(831, 13).Value = 200
$ws.Cells.Item(831, 14).Value = 30000
$ws.Cells.Item(831, 15).Value = 30000
$ws.Cells.Item(831, 16).Value = 30000
$ws.Cells.Item(831, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(831, 18).Value = "Perú"
$ws.Cells.Item(831, 19).Value = 3000
$ws.Cells.Item(831, 20).Value = 10

# Row 832 - Primera
$ws.Cells.Item(832, 1).Value = 11
$ws.Cells.Item(832, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(832, 3).Value = "Bíobío"
$ws.Cells.Item(832, 4).Value = 45021
$ws.Cells.Item(832, 5).Value = 8
$ws.Cells.Item(832, 6).Value = "Fruta"
$ws.Cells.Item(832, 7).Value = 100106
$ws.Cells.Item(832, 8).Value = "Oleaginosos"
$ws.Cells.Item(832, 9).Value = 100106002
$ws.Cells.Item(832, 10).Value = "Palta"
$ws.Cells.Item(832, 11).Value = "Hass"
$ws.Cells.Item(832, 12).Value = "Primera"
$ws.Cells.Item(832, 13).Value = 250
$ws.Cells.Item(832, 14).Value = 28000
$ws.Cells.Item(832, 15).Value = 28000
$ws.Cells.Item(832, 16).Value = 28000
$ws.Cells.Item(832, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(832, 18).Value = "Perú"
$ws.Cells.Item(832, 19).Value = 2800
$ws.Cells.Item(832, 20).Value = 10

# Row 833 - Segunda
$ws.Cells.Item(833, 1).Value = 11
$ws.Cells.Item(833, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(833, 3).Value = "Bíobío"
$ws.Cells.Item(833, 4).Value = 45021
$ws.Cells.Item(833, 5).Value = 8
$ws.Cells.Item(833, 6).Value = "Fruta"
$ws.Cells.Item(833, 7).Value = 100106
$ws.Cells.Item(833, 8).Value = "Oleaginosos"
$ws.Cells.Item(833, 9).Value = 100106002
$ws.Cells.Item(833, 10).Value = "Palta"
$ws.Cells.Item(833, 11).Value = "Hass"
$ws.Cells.Item(833, 12).Value = "Segunda"
$ws.Cells.Item(833, 13).Value = 200
$ws.Cells.Item(833, 14).Value = 25000
$ws.Cells.Item(833, 15).Value = 25000
$ws.Cells.Item(833, 16).Value = 25000
$ws.Cells.Item(833, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(833, 18).Value = "Perú"
$ws.Cells.Item(833, 19).Value = 2500
$ws.Cells.Item(833, 20).Value = 10
